$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated regression results from rerun of lm; permuted models brs, sdnn
$values = @{
    "C2" = 1.65
    "D2" = 0.16
    "E2" = 3.15
    "F2" = 0.03054
    "G2" = 0.03152
    "C3" = -1.1
    "D3" = -2.62
    "E3" = 0.41
    "F3" = 0.15274
    "G3" = 0.52236
    "C4" = 0
    "D4" = -1.42
    "E4" = 1.42
    "F4" = 0.99913
    "G4" = 0.64396
    "C5" = 2.55
    "D5" = 1.09
    "E5" = 4.02
    "F5" = 0.00073
    "G5" = 0.91431
    "C6" = -1.79
    "D6" = -3.33
    "E6" = -0.24
    "F6" = 0.02378
    "G6" = 0.09023
    "C7" = -0.31
    "D7" = -1.8
    "E7" = 1.19
    "F7" = 0.68721
    "G7" = 0.60731
    "C8" = -0.27
    "D8" = -2
    "E8" = 1.46
    "F8" = 0.75797
    "G8" = 0.71809
    "C9" = -0.26
    "D9" = -1.78
    "E9" = 1.26
    "F9" = 0.73797
    "G9" = 0.98574
    "C10" = 0.86
    "D10" = -0.78
    "E10" = 2.5
    "F10" = 0.30416
    "G10" = 0.24896
    "C11" = -0.79
    "D11" = -2.31
    "E11" = 0.73
    "F11" = 0.30754
    "G11" = 0.37328
    "C12" = -0.75
    "D12" = -2.51
    "E12" = 1
    "F12" = 0.39657
    "G12" = 0.03152
    "C13" = -0.75
    "D13" = -2.54
    "E13" = 1.03
    "F13" = 0.40752
    "G13" = 0.52236
    "C14" = 0.71
    "D14" = -0.91
    "E14" = 2.33
    "F14" = 0.38752
    "G14" = 0.64396
    "C15" = 2.34
    "D15" = 0.43
    "E15" = 4.24
    "F15" = 0.01646
    "G15" = 0.91431
    "C16" = 0.23
    "D16" = -1.49
    "E16" = 1.94
    "F16" = 0.79419
    "G16" = 0.09023
    "C17" = -0.05
    "D17" = -1.71
    "E17" = 1.61
    "F17" = 0.95043
    "G17" = 0.60731
    "C18" = 0.5
    "D18" = -1.18
    "E18" = 2.18
    "F18" = 0.55724
    "G18" = 0.71809
    "D19" = -1.86
    "E19" = 1.66
    "F19" = 0.90869
    "G19" = 0.98574
    "C20" = 0.35
    "D20" = -1.27
    "E20" = 1.97
    "F20" = 0.67125
    "G20" = 0.24896
    "C21" = 0.31
    "D21" = -1.41
    "E21" = 2.02
    "F21" = 0.72314
    "G21" = 0.37328
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
